# ZBP_11_obavy_epidemie.xlsx update: add a new survey wave "30. 3. 2021"
# - sheet "data": new column AB with header "30. 3. 2021" and per-row percentage values
# - sheet "pocetR": new column AA with header "30. 3. 2021" and per-row sample-size values
# - update the "aktualizace" date note on both sheets from 23. 3. 2021 to 7. 4. 2021

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data"
$ws2 = $wb.Worksheets.Item(2)   # "pocetR"

# ---- Sheet "data": new column AB ----
$ws1.Range("AB1").Value = "30. 3. 2021"
$ws1.Range("AA1").Copy()
$ws1.Range("AB1").PasteSpecial(-4122)  # xlPasteFormats, match the AA1 header look (bold, border, centered)

$dataValues = @(0.19,0.36,0.45,0.24,0.43,0.33,0.18,0.36,0.46,0.16,0.33,0.51,0.21,0.35,0.44,0.17,0.32,0.51,0.19,0.42,0.39,0.24,0.39,0.37,0.14,0.34,0.52,0.18,0.36,0.46,0.2,0.35,0.45,0.21,0.38,0.41,0.16,0.39,0.45,0.05,0.38,0.57,0.17,0.36,0.47,0.21,0.31,0.48,0.12,0.39,0.49,0.22,0.36,0.42,0.19,0.4,0.41,0.3,0.46,0.24,0.22,0.37,0.41,0.22,0.37,0.41,0.18,0.41,0.41,0.14,0.35,0.51,0.11,0.29,0.6)

for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 28).Value = $dataValues[$i]   # column 28 = AB
}

# ---- Sheet "pocetR": new column AA ----
$ws2.Range("AA1").Value = "30. 3. 2021"
$ws2.Range("Z1").Copy()
$ws2.Range("AA1").PasteSpecial(-4122)  # xlPasteFormats, match the Z1 header look (bold, border, centered)

$countValues = @(2061,488,758,815,626,734,701,1007,1054,1084,464,241,272,51,151,95,20,291,558,251,391,365,238,367,449)

for ($i = 0; $i -lt $countValues.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 27).Value = $countValues[$i]   # column 27 = AA
}

# Row 27 on pocetR has a blank (empty-string-typed) cell under every other data column;
# keep the same pattern under the new AA column.
$ws2.Range("AA27").Value = "trideni"

# ---- Update the "aktualizace" date notes on both sheets ----
$ws1.Range("A77").Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 7. 4. 2021"
$ws2.Range("A27").Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 7. 4. 2021"
